# "remove the top icon flow"
#
# MainIcon sheet has a table (表1, A1:L29) describing icons. The "top icon
# flow" (Flow = 3, rows 17-19: SideButton1/SideButton9/SideButton5) is
# folded into Flow 2, and the two now-orphaned "top flow" extra entries
# (Id 1000 "SideButton10" and Id 1001 "SideButton11") are removed from
# the table entirely. Everything below shifts up to fill the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two obsolete rows (Id 1000 / 1001, SideButton10 / SideButton11).
# Deleting whole rows shifts the rows below (old 25-29 -> new 23-27) up,
# auto-shrinks the table/autofilter range and dimension, and drops the
# now-unused shared strings.
$ws.Rows("23:24").Delete()

# Re-point the old "top flow" (3) rows at flow 2, merging the flows.
$ws.Range("K17:K19").Value = 2

# Leave the selection where the edit ended up.
$ws.Range("K17:K27").Select()
